# Append a new attendance date column (U) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell U1: new date column header.
# Write the date as literal text (leading quote) first, so Excel doesn't
# reinterpret the "2025-06-23" string as a date serial value, then copy
# just the formatting (bold, centered, bordered) over from T1 so U1 matches
# the look of the rest of the header row.
$ws.Range("U1").Value = "'2025-06-23"
$ws.Range("T1").Copy() | Out-Null
$ws.Range("U1").PasteSpecial(-4122)

# Row 2: mark absent for the new date, and bump the running Total (S2) by 1.
$ws.Range("U2").Value = "❌"
$ws.Range("S2").Value = 16

# Row 3: same as row 2.
$ws.Range("U3").Value = "❌"
$ws.Range("S3").Value = 16
